$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9719789624214172
$ws.Range("B1").Value = 0.7558236122131348
$ws.Range("C1").Value = 3.368383646011353
$ws.Range("D1").Value = 3.217042684555054
$ws.Range("E1").Value = 0.9465873837471008
